# Apply the "phase3" budget-sheet update:
#  - rename TRICALL to Triangulation add RefinePart, EdgeLengths add DepthFirstFindAllFaces
#  - highlight several finished tasks in column B with the built-in "Good" style
#  - add two inline remarks in column I (rows 7 and 17)
#  - add a new "hours spent" mini-table in rows 27-36 with totals

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column A slightly to fit the new "hours spent" label
$ws.Columns.Item(1).ColumnWidth = 25.59

# Mark these completed/important tasks with the built-in "Good" (green) style
$goodCells = @("B7", "B10", "B11", "B12", "B14", "B15", "B16", "B17", "B18")
foreach ($cellRef in $goodCells) {
    $ws.Range($cellRef).Style = "Good"
}

# Inline remarks
$ws.Range("I7").Value = "make interactive, solve last bugs"
$ws.Range("I17").Value = "this should be easy"

# Replace the old stray formatted-but-empty B27 cell with the new
# "hours spent" mini table (rows 27-34) and totals row (36).
# (C29 keeps its existing style untouched.)
$ws.Range("B27").ClearFormats()

$ws.Range("A27").Value = "hours spent"
$ws.Range("B27").Value = 22
$ws.Range("B28").Value = 20
$ws.Range("B29").Value = 26
$ws.Range("B30").Value = 20
$ws.Range("B31").Value = 31
$ws.Range("B32").Value = 20
$ws.Range("B33").Value = 36
$ws.Range("B34").Value = 36

$ws.Range("A36").Value = "Total hours spent so far"
$ws.Range("B36").Formula = "=SUM(B27:B34)"
$ws.Range("C36").Formula = "=B36*135"

$ws.Range("A36").Select() | Out-Null
